{"js": "// Work item iz\u00e9k initial\n//\n// 1) \"Megrendel\u0151: \" -> \"Megrendel\u0151: asdasd\"\n// 2) \"2024.06.08\" -> \"2024.06.10\" (three occurrences across the document)\n\n// --- Change 1: append \"asdasd\" to the \"Megrendel\u0151: \" label -----------------\nconst megrendeloResults = context.document.body.search(\"Megrendel\u0151: \", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nmegrendeloResults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < megrendeloResults.items.length; i++) {\n  const item = megrendeloResults.items[i];\n  item.insertText(item.text + \"asdasd\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Change 2: bump every \"2024.06.08\" date to \"2024.06.10\" ----------------\nconst dateResults = context.document.body.search(\"2024.06.08\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\ndateResults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  const item = dateResults.items[i];\n  const updatedText = item.text.split(\"2024.06.08\").join(\"2024.06.10\");\n  item.insertText(updatedText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Work item iz\u00e9k initial\n#\n# 1) \"Megrendel\u0151: \" -> \"Megrendel\u0151: asdasd\"\n# 2) \"2024.06.08\" -> \"2024.06.10\" (three occurrences across the document)\n\n$d = $word.ActiveDocument\n\n# --- Change 1: append \"asdasd\" to the \"Megrendel\u0151: \" label -----------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Megrendel\u0151: \"\n$find1.MatchCase = $true\n$find1.Replacement.Text = \"Megrendel\u0151: asdasd\"\n$find1.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# --- Change 2: bump every \"2024.06.08\" date to \"2024.06.10\" ----------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"2024.06.08\"\n$find2.MatchCase = $true\n$find2.Replacement.Text = \"2024.06.10\"\n$find2.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n"}
